$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# RNF007 row (row 9): "Biografia(Perfil)" requirement replaced by an avatar requirement.
$ws.Range("D9").Value = "Foto de Perfil"
$ws.Range("E9").Value = "O sistema deve exibir o avatar padrão do usuário na barra lateral."

# RF012 row (row 22): used to be the "Perfil" listing requirement, now the Dashboard stats requirement.
$ws.Range("D22").Value = "Dashboard"
$ws.Range("E22").Value = "O sistema deve exibir estatísticas como número de posts, curtidas e comentários recebidos."

# RF013 row (row 23): keeps "Dashboard" but the requirement text becomes the engagement-chart one.
$ws.Range("E23").Value = "O sistema deve exibir gráficos simples com engajamento ao longo do tempo."

# The old RF014 ("gráficos simples...") and RF015 ("sugerir ações...") rows are removed entirely.
$ws.Range("A24:A25").EntireRow.Delete()

# Update the active selection / view to match where the edits were made.
$ws.Range("F23").Select()
